$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A3").Value = 60475415
$ws.Range("S3").Value = 5

$ws.Range("Y3").NumberFormat = "@"
$ws.Range("Y3").Value = "2016-07-06"

$ws.Range("AA3").NumberFormat = "@"
$ws.Range("AA3").Value = "2016-07-06"

$ws.Range("AC3").Value = "Enstaka blommande ex"
$ws.Range("AW3").Value = "Göran Frisk"
$ws.Range("AX3").Value = "Göran Frisk, Kristina Nygren Frisk"
